$d = $word.ActiveDocument

# The document starts with a paragraph containing a manual page break
# (<w:r><w:br w:type="page"/></w:r>) followed by a section-break paragraph
# and a final paragraph that only holds the hidden "_GoBack" bookmark.
#
# The edit removes the page break character from the first paragraph and
# relocates the "_GoBack" bookmark from the last paragraph to that now
# empty first paragraph, leaving the last paragraph completely empty.

# 1) Delete the page-break character (the single character run) that is
#    the entire content of the first paragraph.
$d.Range(0, 1).Delete()

# 2) Remove the existing "_GoBack" bookmark (currently around the last
#    paragraph) so it can be re-created at the start of the document.
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 3) Re-insert the "_GoBack" bookmark, collapsed, at the very start of the
#    document - i.e. inside the first paragraph, which is now empty.
$d.Bookmarks.Add("_GoBack", $d.Range(0, 0))
